# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the newly generated output data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Changes for sheet "展览" (row => new value)
$sheet1Changes = @{
    2  = 599
    4  = 1277
    5  = 1132
    6  = 14199
    7  = 15973
    8  = 14
    9  = 68
    10 = 47
    20 = 1234
    23 = 25
    24 = 6349
    25 = 964
    27 = 5632
    28 = 82
    30 = 141
    31 = 4630
}

foreach ($row in $sheet1Changes.Keys) {
    $ws1.Range("F$row").Value = $sheet1Changes[$row]
}

# Changes for sheet "全部类型" (row => new value)
$sheet4Changes = @{
    2  = 599
    4  = 1277
    5  = 1132
    6  = 14199
    7  = 15973
    8  = 14
    9  = 68
    10 = 47
    20 = 1234
    24 = 25
    25 = 6349
    26 = 964
    29 = 5632
    30 = 82
    32 = 141
    33 = 4630
}

foreach ($row in $sheet4Changes.Keys) {
    $ws4.Range("F$row").Value = $sheet4Changes[$row]
}
